$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header column D: alerttext
$ws.Range("D1").Value = "alerttext"

# New alert text value added first (so shared-string order matches: alerttext, Customer added successfully, Rakesh1, Bhavsar1)
$ws.Range("D2").Value = "Customer added successfully"

# Update existing row 2 values (Rakesh -> Rakesh1, Bhavsar -> Bhavsar1, postCode -> 4240011)
$ws.Range("A2").Value = "Rakesh1"
$ws.Range("B2").Value = "Bhavsar1"
$ws.Range("C2").Value = 4240011

# Update selection to match target sheet view
$ws.Range("D5").Select()
